$d = $word.ActiveDocument

# 1) The "OPTIONAL Inclusion of | . | ! | @ | # | ..." bullet becomes "Case sensitive".
$pOptional = $d.Paragraphs(3)
$rOptional = $pOptional.Range
$rOptional.MoveEnd(1, -1) | Out-Null
$rOptional.Text = "Case sensitive"

# 2) The old "Case sensitive" bullet becomes "A-Z and a-z and 0-9".
$pCaseSensitive = $d.Paragraphs(4)
$rCaseSensitive = $pCaseSensitive.Range
$rCaseSensitive.MoveEnd(1, -1) | Out-Null
$rCaseSensitive.Text = "A-Z and a-z and 0-9"

# 3) The old "A-Z and a-z and 0-9" bullet is removed entirely (whole paragraph, incl. mark).
$pOldAZ = $d.Paragraphs(5)
$pOldAZ.Range.Delete() | Out-Null

# 4) The "_GoBack" bookmark moves from the end of "No spaces allowed" to the very
#    start of the (now renumbered) "Case sensitive" bullet.
$d.Bookmarks("_GoBack").Delete()
$pNowCaseSensitive = $d.Paragraphs(3)
$rStart = $pNowCaseSensitive.Range
$rStart.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $rStart) | Out-Null
